$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: rowNumber, startColumnIndex(1-based), then contiguous new values for that row run
$rowUpdates = @(
    ,@(8, 4, 72158800, 66387100, 68941800, 69743900, 69938300, 66019500, 70929800)
    ,@(9, 4, 51012300, 46620800, 96797300, 49966000, 50975400, 48995800, 53015200)
    ,@(10, 4, 21146500, 19766400, -27855500, 19777900, 18963000, 17023700, 17914600)
    ,@(12, 6, 4066400, 4133500, 4328500, 4540100, 4702800)
    ,@(14, 4, 51500, 191700, 1271100, 506500, 1011000, 3514000, 5091000)
    ,@(15, 4, 767900, 745300, 2890500, 2189000, 2520300, 2509300, 2342600)
    ,@(17, 4, 68718700, 63885000, 66859900, 66798000, 68191100, 68078700, 75625500)
    ,@(18, 4, 3440100, 2502100, 2081900, 2946000, 1747300, -2059200, -4695700)
    ,@(20, 4, 205900, 197400, 367300, -1137800, 315100, -1310800, -2395600)
    ,@(21, 4, 6250100, 5149900, 4968600, 4401300, 5058700, -298700, -4031400)
    ,@(22, 4, 223500, 212900, 392300, 158800, 198100, 231400, 256800)
    ,@(23, 4, 3422500, 2486600, 2056900, 1649400, 1864300, -3601400, -7348100)
    ,@(24, 4, 1144100, 927700, 328100, -17900, 810600, 3477400, 88300)
    ,@(26, 4, 2278300, 1558900, 1728700, 1667300, 1053700, -7078900, -7436400)
    ,@(27, 4, 2133800, 1350200, 1409900, 1622500, 1088800, -6818400, -6980400)
    ,@(32, 4, -205900, -197400, -367300, 1137800, -315100, 1310800, 2395600)
    ,@(33, 4, 2133800, 1350200, 1409900, 1622500, 1088800, -6818400, -6980400)
    ,@(35, 4, 2133800, 1350200, 1409900, 1622500, 1088800, -6818400, -6980400)
    ,@(41, 4, 9849800, 11487900, 18324800, 11741900, 5355900, 4501500, 5523300)
    ,@(42, 4, 1840200, 1297400, 1496100)
    ,@(42, 10, 4400)
    ,@(43, 4, 9392400, 7656900, 14997400, 8968500, 9107200, 8491500, 9127200)
    ,@(44, 4, 8937000, 7289000, 13795900, 6894500, 6786200, 7113100, 7505600)
    ,@(45, 4, 1493600, 1240300, 5156800, 3246200, 2742800, 2440400, 4110200)
    ,@(46, 4, 31513100, 28971600, 26159600, 30851200, 23992100, 22546500, 26270600)
    ,@(47, 4, 2842400, 2874500, 5917500, 2835600, 2457100, 2503900, 4085000)
    ,@(48, 4, 12421600, 11962500, 23408300, 12428500, 12886100, 15145900, 15677900)
    ,@(49, 4, 6673800, 6012800, 9870200, 5695200, 6102600, 6645800, 9972600)
    ,@(52, 4, 3421200, 4264700, 7241500, 2040400, 1687600, 1954100, 3667400)
    ,@(54, 4, 56872000, 54086000, 49611700, 53850800, 47125500, 48796200, 59673500)
    ,@(57, 4, 10364100, 8641900, 16608000, 8889100, 8470600, 7157700, 7693200)
    ,@(58, 4, 3393500, 1600400, 392800, 2355200, 766000, 4341900, 5730000)
    ,@(59, 4, 14247600, 14274700, 26971900, 13460200, 12801600, 11996700, 12607600)
    ,@(60, 4, 28005300, 24517000, 22449400, 24704500, 22038200, 23496400, 26030700)
    ,@(61, 4, 7811000, 8560600, 6356100, 6440000, 5038700, 5994300, 8513600)
    ,@(62, 4, 4039800, 5098500, 11859700, 4693700, 5707200, 7514900, 7252000)
    ,@(66, 4, 41435700, 39876100, 36554000, 37368200, 33130200, 37369400, 42228300)
    ,@(72, 4, 11755000, 9505100, 18473100, 9232000, 7943800, 6959600, 13882600)
    ,@(76, 4, 15436300, 14209900, 13057800, 16482600, 13995300, 11426800, 17445300)
    ,@(81, 4, 2133800, 1350200, 1409900, 1622500, 1088800, -6818400, -6980400)
    ,@(83, 4, 2601300, 2447700, 2516700, 2590200, 2993000, 3067900, 3056500)
    ,@(89, 4, 3825600, 3484100, 3791000, 4442800, 5260800, 3062300, -333500)
    ,@(91, 4, -3566100, -2518500, -4403400, -2026400, -1823700, -2894300, -4126500)
    ,@(94, 4, -4147800, -3798200, -2656000, -1247600, 109600, 148300, -2739100)
    ,@(96, 4, -527100, -524500, -418800, -334300, -104500, -104500, -198100)
    ,@(100, 4, -1164000, 2663200, -2798500, 2328800, -4812100, -4439200, -480000)
    ,@(101, 4, -151800, -15600, -752700, 694900, 311200, 522300, -67100)
    ,@(102, 4, -1638100, 2333400, -2416200, 6219000, 869500, -706300, -3619800)
)

foreach ($entry in $rowUpdates) {
    $rowNum = $entry[0]
    $startCol = $entry[1]
    for ($i = 2; $i -lt $entry.Length; $i++) {
        $ws.Cells.Item($rowNum, $startCol + $i - 2).Value2 = $entry[$i]
    }
}
